$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: the sentence about ledgers being transparent/unchangeable
# was fragmented across several <w:r> elements in the source markup.
# The visible text itself is unchanged; running it through
# Find/Replace consolidates the fragmented runs into a single run.
# -----------------------------------------------------------------
$ledgersFind = $d.Content.Find
$ledgersFind.ClearFormatting()
$ledgersFind.Replacement.ClearFormatting()
$ledgersText = ", these ledgers are transparent, unchangeable, and stored as verified copy" + [char]8217 + "s between multiple nodes."
$ledgersFind.Execute($ledgersText, $true, $false, $false, $false, $false, $true, 1, $false, $ledgersText, 2) | Out-Null

# -----------------------------------------------------------------
# Change 2: fix capitalisation and remove the duplicate paragraph
# break between "...essentially the same thing." and "these mediums
# of exchange...". The two paragraphs are merged into one, and
# "these" becomes " These" (space inserted, "t" capitalised to "T").
# -----------------------------------------------------------------

# Find the paragraph that currently ends with "...essentially the same
# thing." and delete its trailing paragraph mark, merging it with the
# following paragraph ("these mediums of exchange...").
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "essentially the same thing\.\s*$") {
        $markRange = $d.Range($p.Range.End - 1, $p.Range.End)
        if ($markRange.Text -eq [string][char]13) {
            $markRange.Delete()
        }
        break
    }
}

# Now that the paragraphs are merged, fix up the capitalisation in one
# pass: "same thing.these mediums" -> "same thing. These mediums".
$capFind = $d.Content.Find
$capFind.ClearFormatting()
$capFind.Replacement.ClearFormatting()
$capFind.Execute("same thing.these mediums", $true, $false, $false, $false, $false, $true, 1, $false, "same thing. These mediums", 2) | Out-Null
